# Fruta / hortaliza, semanal
#
# Inserts a new weekly price record as a new row 364 in the "Hortaliza,
# Feria Lagunitas de Puerto Montt - Ciboulette" sheet, shifting the
# existing rows 364:398 down to 365:399.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 364; everything currently
# at row 364 (and below, through 398) moves down to 365 (through 399).
$ws.Rows.Item(364).Insert()

# Populate the newly inserted row 364 with the new weekly record.
$ws.Cells.Item(364, 1).Value = 4
$ws.Cells.Item(364, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(364, 3).Value = "Los Lagos"
$ws.Cells.Item(364, 4).Value = 45166
$ws.Cells.Item(364, 5).Value = 10
$ws.Cells.Item(364, 6).Value = 100112039
$ws.Cells.Item(364, 7).Value = "Ciboulette"
$ws.Cells.Item(364, 8).Value = "Sin especificar"
$ws.Cells.Item(364, 9).Value = "Primera"
$ws.Cells.Item(364, 10).Value = 80
$ws.Cells.Item(364, 11).Value = 3500
$ws.Cells.Item(364, 12).Value = 3500
$ws.Cells.Item(364, 13).Value = 3500
$ws.Cells.Item(364, 14).Value = "$/docena de atados"
$ws.Cells.Item(364, 15).Value = "Región Metropolitana"
$ws.Cells.Item(364, 16).Value = 1167
$ws.Cells.Item(364, 17).Value = 3
$ws.Cells.Item(364, 18).Value = "Hortaliza"
